$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 691.28986
$ws.Range("J17").Value = 585.2295
$ws.Range("L17").Value = 1755.6885
$ws.Range("N17").Value = -2091.6885
# Row 19
$ws.Range("H19").Value = 1549087.6
$ws.Range("I19").Value = 2632327
$ws.Range("J19").Value = 1602.7142
$ws.Range("K19").Value = 2632327
$ws.Range("L19").Value = 1602.7142
$ws.Range("M19").Value = -2632152
$ws.Range("N19").Value = -1952.7142
# Row 100
$ws.Range("H100").Value = 11766382
$ws.Range("I100").Value = 12501469
$ws.Range("K100").Value = 12501469
$ws.Range("M100").Value = -12500928
# Row 106
$ws.Range("H106").Value = 4600.4546
$ws.Range("I106").Value = 1967.5
$ws.Range("K106").Value = 1967.5
$ws.Range("M106").Value = -1336.5
# Row 112
$ws.Range("H112").Value = 27028270
$ws.Range("I112").Value = 333333800
$ws.Range("J112").Value = 1312.7354
$ws.Range("K112").Value = 1000001400
$ws.Range("L112").Value = 3938.2062
$ws.Range("M112").Value = -1000000292
$ws.Range("N112").Value = -6154.206200000001
# Row 113
$ws.Range("H113").Value = 4285.4287
$ws.Range("J113").Value = 4750.9165
$ws.Range("L113").Value = 4750.9165
$ws.Range("N113").Value = -11258.9165
# Row 116
$ws.Range("H116").Value = 424121.25
$ws.Range("I116").Value = 912100.4399999999
$ws.Range("J116").Value = 11215.77
$ws.Range("K116").Value = 912100.4399999999
$ws.Range("L116").Value = 11215.77
$ws.Range("M116").Value = -908658.4399999999
$ws.Range("N116").Value = -18099.77
# Row 137
$ws.Range("H137").Value = 1362669.2
$ws.Range("I137").Value = 2647164.2
$ws.Range("K137").Value = 7941492.600000001
$ws.Range("M137").Value = -7938942.600000001
# Row 138
$ws.Range("H138").Value = 4792.35
$ws.Range("I138").Value = 926.3333
$ws.Range("J138").Value = 6013.1973
$ws.Range("K138").Value = 2778.9999
$ws.Range("L138").Value = 18039.5919
$ws.Range("M138").Value = 2361.0001
$ws.Range("N138").Value = -28319.5919

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3577.8293
$ws.Range("I61").Value = 1087.875
$ws.Range("J61").Value = 12431
$ws.Range("K61").Value = 1087.875
$ws.Range("L61").Value = 12431
$ws.Range("M61").Value = -875.875
$ws.Range("N61").Value = -12855
# Row 74
$ws.Range("H74").Value = 5049.231
$ws.Range("I74").Value = 7143.769
$ws.Range("J74").Value = 2954.6924
$ws.Range("K74").Value = 7143.769
$ws.Range("L74").Value = 2954.6924
$ws.Range("M74").Value = -6269.769
$ws.Range("N74").Value = -4702.6924
# Row 77
$ws.Range("H77").Value = 5049.231
$ws.Range("I77").Value = 7143.769
$ws.Range("J77").Value = 2954.6924
$ws.Range("K77").Value = 35718.845
$ws.Range("L77").Value = 14773.462
$ws.Range("M77").Value = -31350.845
$ws.Range("N77").Value = -23509.462
# Row 132
$ws.Range("H132").Value = 1431.88
$ws.Range("I132").Value = 908.64
$ws.Range("J132").Value = 2478.36
$ws.Range("K132").Value = 2725.92
$ws.Range("L132").Value = 7435.08
$ws.Range("M132").Value = -195.9200000000001
$ws.Range("N132").Value = -12495.08
# Row 136
$ws.Range("H136").Value = 3577.8293
$ws.Range("I136").Value = 1087.875
$ws.Range("J136").Value = 12431
$ws.Range("K136").Value = 3263.625
$ws.Range("L136").Value = 37293
$ws.Range("M136").Value = -713.625
$ws.Range("N136").Value = -42393

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 4097.857
$ws.Range("I99").Value = 1497.5
$ws.Range("J99").Value = 5138
$ws.Range("K99").Value = 1497.5
$ws.Range("L99").Value = 5138
$ws.Range("M99").Value = 0.5
$ws.Range("N99").Value = -8134
# Row 134
$ws.Range("H134").Value = 4549.317
$ws.Range("I134").Value = 1466.742
$ws.Range("K134").Value = 4400.226
$ws.Range("M134").Value = -1865.226

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2584.8684
$ws.Range("I58").Value = 1601.2812
$ws.Range("J58").Value = 7830.6665
$ws.Range("K58").Value = 1601.2812
$ws.Range("L58").Value = 7830.6665
$ws.Range("M58").Value = -1398.2812
$ws.Range("N58").Value = -8236.666499999999
# Row 132
$ws.Range("H132").Value = 2825.122
$ws.Range("I132").Value = 2436.2354
$ws.Range("K132").Value = 7308.706200000001
$ws.Range("M132").Value = -4778.706200000001
# Row 134
$ws.Range("H134").Value = 2274.261
$ws.Range("I134").Value = 1312.2727
$ws.Range("J134").Value = 3156.0833
$ws.Range("K134").Value = 3936.8181
$ws.Range("L134").Value = 9468.249899999999
$ws.Range("M134").Value = -1401.8181
$ws.Range("N134").Value = -14538.2499
# Row 136
$ws.Range("H136").Value = 2584.8684
$ws.Range("I136").Value = 1601.2812
$ws.Range("J136").Value = 7830.6665
$ws.Range("K136").Value = 4803.8436
$ws.Range("L136").Value = 23491.9995
$ws.Range("M136").Value = -2253.8436
$ws.Range("N136").Value = -28591.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 37993.75
$ws.Range("I4").Value = 100083.336
$ws.Range("J4").Value = 740
$ws.Range("K4").Value = 300250.008
$ws.Range("L4").Value = 2220
$ws.Range("M4").Value = -300138.008
$ws.Range("N4").Value = -2444
# Row 12
$ws.Range("H12").Value = 92.26667
$ws.Range("J12").Value = 112.75
$ws.Range("L12").Value = 338.25
$ws.Range("N12").Value = -684.25
# Row 70
$ws.Range("H70").Value = 2245
$ws.Range("I70").Value = 1360.1428
$ws.Range("K70").Value = 4080.4284
$ws.Range("M70").Value = -3765.4284
# Row 73
$ws.Range("H73").Value = 2245
$ws.Range("I73").Value = 1360.1428
$ws.Range("K73").Value = 4080.4284
$ws.Range("M73").Value = -2988.4284
# Row 92
$ws.Range("H92").Value = 296.66666
$ws.Range("I92").Value = 296.66666
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 889.9999799999999
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 358.0000200000001
$ws.Range("N92").ClearContents()
# Row 127
$ws.Range("H127").Value = 916.5
$ws.Range("J127").Value = 916.5
$ws.Range("L127").Value = 2749.5
$ws.Range("N127").Value = -12669.5
# Row 129
$ws.Range("H129").Value = 2416.25
$ws.Range("J129").Value = 2641
$ws.Range("L129").Value = 7923
$ws.Range("N129").Value = -17923
# Row 131
$ws.Range("H131").Value = 777.63
$ws.Range("I131").Value = 387.16666
$ws.Range("J131").Value = 802.55316
$ws.Range("K131").Value = 1161.49998
$ws.Range("L131").Value = 2407.65948
$ws.Range("M131").Value = 3878.50002
$ws.Range("N131").Value = -12487.65948

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2177.913
$ws.Range("I102").Value = 1594.6
$ws.Range("J102").Value = 6066.6665
$ws.Range("K102").Value = 1594.6
$ws.Range("L102").Value = 6066.6665
$ws.Range("M102").Value = 27.40000000000009
$ws.Range("N102").Value = -9310.666499999999
# Row 132
$ws.Range("H132").Value = 2215.8386
$ws.Range("I132").Value = 1429.1818
$ws.Range("K132").Value = 4287.5454
$ws.Range("M132").Value = -1757.5454

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4789.25
$ws.Range("I7").Value = 3680.7
$ws.Range("J7").Value = 6636.8335
$ws.Range("K7").Value = 3680.7
$ws.Range("L7").Value = 6636.8335
$ws.Range("M7").Value = -3568.7
$ws.Range("N7").Value = -6860.8335
# Row 11
$ws.Range("H11").Value = 7004751
$ws.Range("J11").Value = 7004751
$ws.Range("L11").Value = 7004751
$ws.Range("N11").Value = -7005031
# Row 93
$ws.Range("H93").Value = 4117501.5
$ws.Range("I93").Value = 8549271
$ws.Range("J93").Value = 2287.0715
$ws.Range("K93").Value = 8549271
$ws.Range("L93").Value = 2287.0715
$ws.Range("M93").Value = -8548023
$ws.Range("N93").Value = -4783.0715
# Row 126
$ws.Range("H126").Value = 4789.25
$ws.Range("I126").Value = 3680.7
$ws.Range("J126").Value = 6636.8335
$ws.Range("K126").Value = 11042.1
$ws.Range("L126").Value = 19910.5005
$ws.Range("M126").Value = -8572.099999999999
$ws.Range("N126").Value = -24850.5005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2587.818
$ws.Range("I136").Value = 1704.909
$ws.Range("K136").Value = 5114.727000000001
$ws.Range("M136").Value = -2564.727000000001
